# Auto-generated Excel COM-interop script
# Applies numeric corrections to currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1112.75
$ws.Range("J19").Value = 1200.5
$ws.Range("L19").Value = 1200.5
$ws.Range("N19").Value = -1550.5
$ws.Range("H115").Value = 631.4286
$ws.Range("I115").Value = 631.4286
$ws.Range("K115").Value = 1894.2858
$ws.Range("M115").Value = -327.2857999999999
$ws.Range("H135").Value = 2108.5217
$ws.Range("I135").Value = 2024.8
$ws.Range("K135").Value = 18223.2
$ws.Range("M135").Value = -15688.2
$ws.Range("H138").Value = 2772.1194
$ws.Range("I138").Value = 1477.8928
$ws.Range("J138").Value = 3701.3076
$ws.Range("K138").Value = 4433.678400000001
$ws.Range("L138").Value = 11103.9228
$ws.Range("M138").Value = 706.3215999999993
$ws.Range("N138").Value = -21383.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2997.3333
$ws.Range("I61").Value = 2935.125
$ws.Range("J61").Value = 3495
$ws.Range("K61").Value = 2935.125
$ws.Range("L61").Value = 3495
$ws.Range("M61").Value = -2723.125
$ws.Range("N61").Value = -3919
$ws.Range("H74").Value = 1197.3334
$ws.Range("I74").Value = 1224.0212
$ws.Range("J74").Value = 1071.9
$ws.Range("K74").Value = 1224.0212
$ws.Range("L74").Value = 1071.9
$ws.Range("M74").Value = -350.0211999999999
$ws.Range("N74").Value = -2819.9
$ws.Range("H77").Value = 1197.3334
$ws.Range("I77").Value = 1224.0212
$ws.Range("J77").Value = 1071.9
$ws.Range("K77").Value = 6120.106
$ws.Range("L77").Value = 5359.5
$ws.Range("M77").Value = -1752.106
$ws.Range("N77").Value = -14095.5
$ws.Range("H86").Value = 59995
$ws.Range("J86").Value = 59995
$ws.Range("L86").Value = 59995
$ws.Range("N86").Value = -62367
$ws.Range("H89").Value = 59995
$ws.Range("J89").Value = 59995
$ws.Range("L89").Value = 179985
$ws.Range("N89").Value = -191841
$ws.Range("H124").Value = 24992.5
$ws.Range("J124").Value = 24992.5
$ws.Range("L124").Value = 24992.5
$ws.Range("N124").Value = -34812.5
$ws.Range("H130").Value = 189998.33
$ws.Range("J130").Value = 189998.33
$ws.Range("L130").Value = 189998.33
$ws.Range("N130").Value = -200038.33
$ws.Range("H132").Value = 2075.5386
$ws.Range("I132").Value = 2075.5386
$ws.Range("K132").Value = 6226.6158
$ws.Range("M132").Value = -3696.6158
$ws.Range("H136").Value = 2997.3333
$ws.Range("I136").Value = 2935.125
$ws.Range("J136").Value = 3495
$ws.Range("K136").Value = 8805.375
$ws.Range("L136").Value = 10485
$ws.Range("M136").Value = -6255.375
$ws.Range("N136").Value = -15585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 32880.5
$ws.Range("J119").Value = 32880.5
$ws.Range("L119").Value = 32880.5
$ws.Range("N119").Value = -42556.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2114.9565
$ws.Range("I58").Value = 1535.7778
$ws.Range("K58").Value = 1535.7778
$ws.Range("M58").Value = -1332.7778
$ws.Range("H92").Value = 24866.834
$ws.Range("J92").Value = 24866.834
$ws.Range("L92").Value = 24866.834
$ws.Range("N92").Value = -29858.834
$ws.Range("H132").Value = 3740.7058
$ws.Range("I132").Value = 3191.0908
$ws.Range("J132").Value = 4748.3335
$ws.Range("K132").Value = 9573.2724
$ws.Range("L132").Value = 14245.0005
$ws.Range("M132").Value = -7043.2724
$ws.Range("N132").Value = -19305.0005
$ws.Range("H136").Value = 2114.9565
$ws.Range("I136").Value = 1535.7778
$ws.Range("K136").Value = 4607.3334
$ws.Range("M136").Value = -2057.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1538.1177
$ws.Range("J5").Value = 1900
$ws.Range("L5").Value = 5700
$ws.Range("N5").Value = -5924
$ws.Range("H34").Value = 994.7
$ws.Range("J34").Value = 1374.25
$ws.Range("L34").Value = 4122.75
$ws.Range("N34").Value = -4290.75
$ws.Range("H113").Value = 1457.75
$ws.Range("J113").Value = 1605.75
$ws.Range("L113").Value = 4817.25
$ws.Range("N113").Value = -9157.25
$ws.Range("H134").Value = 3179.5715
$ws.Range("I134").Value = 1962.6154
$ws.Range("K134").Value = 5887.8462
$ws.Range("M134").Value = -817.8462
$ws.Range("H135").Value = 1538.1177
$ws.Range("J135").Value = 1900
$ws.Range("L135").Value = 17100
$ws.Range("N135").Value = -22170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4424800
$ws.Range("I11").Value = 5500000
$ws.Range("K11").Value = 5500000
$ws.Range("M11").Value = -5499861
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H86").Value = 60856.832
$ws.Range("J86").Value = 60856.832
$ws.Range("L86").Value = 60856.832
$ws.Range("N86").Value = -63228.832
$ws.Range("H89").Value = 60856.832
$ws.Range("J89").Value = 60856.832
$ws.Range("L89").Value = 182570.496
$ws.Range("N89").Value = -194426.496
$ws.Range("M12").ClearContents()
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H22").Value = 1316.4166
$ws.Range("I22").Value = 849.75
$ws.Range("J22").Value = 1549.75
$ws.Range("K22").Value = 849.75
$ws.Range("L22").Value = 1549.75
$ws.Range("M22").Value = -554.75
$ws.Range("N22").Value = -2139.75
$ws.Range("H27").Value = 1316.4166
$ws.Range("I27").Value = 849.75
$ws.Range("J27").Value = 1549.75
$ws.Range("K27").Value = 849.75
$ws.Range("L27").Value = 1549.75
$ws.Range("M27").Value = -742.75
$ws.Range("N27").Value = -1763.75
$ws.Range("H88").Value = 27166.666
$ws.Range("I88").Value = 15000
$ws.Range("J88").Value = 33250
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 33250
$ws.Range("M88").Value = -14572
$ws.Range("N88").Value = -34106
$ws.Range("H91").Value = 27166.666
$ws.Range("I91").Value = 15000
$ws.Range("J91").Value = 33250
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 33250
$ws.Range("M91").Value = -13518
$ws.Range("N91").Value = -36214
$ws.Range("H109").Value = 41250
$ws.Range("J109").Value = 41250
$ws.Range("L109").Value = 41250
$ws.Range("N109").Value = -44024
$ws.Range("H136").Value = 2151.9092
$ws.Range("I136").Value = 2109.7812
$ws.Range("K136").Value = 6329.3436
$ws.Range("M136").Value = -3779.3436
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 36727.273
$ws.Range("J109").Value = 36727.273
$ws.Range("L109").Value = 36727.273
$ws.Range("N109").Value = -39501.273
$ws.Range("H122").Value = 2115.8667
$ws.Range("I122").Value = 2524.3
$ws.Range("J122").Value = 1299
$ws.Range("K122").Value = 7572.900000000001
$ws.Range("L122").Value = 3897
$ws.Range("M122").Value = -5122.900000000001
$ws.Range("N122").Value = -8797
$ws.Range("H136").Value = 923.9259
$ws.Range("I136").Value = 863.3077
$ws.Range("K136").Value = 2589.9231
$ws.Range("M136").Value = -39.92309999999998
